$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.184.20"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6590"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07429"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.96"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.912.72"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.983"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6657"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("E15").Value = "  -3.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.114"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008612"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.205.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.123.86"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "226.93"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.117"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.604"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1397"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.115"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.046"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05262"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7382"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.146"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.657"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.305.48"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01796"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.733"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9195"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.046"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("B43").Value = "XinFinNetwork"
$ws.Range("C43").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08694"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +13.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.51"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.024.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5142"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.752"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05845"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.04%  "
